$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J3").Value = "BE_001, BE_002"
$ws.Range("J7").Select()
